$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the worker row for HAINER JOSE AHUMADA BATISTA (1050975677) - row 19.
# Deleting the whole row shifts rows 20-25 up by one (so old rows 24/25 become 23/24).
$ws.Rows("19:19").Delete()

# The "Periodo Mora" value shown for all remaining workers changes from 2508 to 2509.
$ws.Range("E16:E18").Value = "2509"

# Center-align the "Periodo Mora" column for the remaining data rows.
$ws.Range("E16:E18").HorizontalAlignment = -4108

# Update the total "VALOR MORA" figure.
$ws.Range("E11").Value = 170820

# Update the worker count ("Cant. Trabajadores") now that one worker was removed.
$ws.Range("C13").Value = 3
